# Tele2 new web page setup
# Append 49 new phone-model mapping rows (rows 236-284) to the "map" sheet,
# extending the used range from A1:D235 to A1:D284.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A236").Value = "Samsung Galaxy S20 Ultra Dual SIM Cosmic Black"
$ws.Range("B236").Value = "SAMSUNG"
$ws.Range("C236").Value = "GALAXY S20 ULTRA"
$ws.Range("D236").Value = "INFONA"

$ws.Range("A237").Value = "`niPhone 11 64GB  `n"
$ws.Range("B237").Value = "APPLE"
$ws.Range("C237").Value = "IPHONE 11"
$ws.Range("D237").Value = "64GB"

$ws.Range("A238").Value = "`nIphone 11 Pro 64GB`n"
$ws.Range("B238").Value = "APPLE"
$ws.Range("C238").Value = "IPHONE 11 PRO"
$ws.Range("D238").Value = "64GB"

$ws.Range("A239").Value = "`niPhone 7 32 GB`n"
$ws.Range("B239").Value = "APPLE"
$ws.Range("C239").Value = "IPHONE 7"
$ws.Range("D239").Value = "32GB"

$ws.Range("A240").Value = "`niPhone 8 64 GB`n"
$ws.Range("B240").Value = "APPLE"
$ws.Range("C240").Value = "IPHONE 8"
$ws.Range("D240").Value = "64GB"

$ws.Range("A241").Value = "`niPhone XR 64 GB`n"
$ws.Range("B241").Value = "APPLE"
$ws.Range("C241").Value = "IPHONE XR"
$ws.Range("D241").Value = "64GB"

$ws.Range("A242").Value = "`niPhone XS 64 GB `n"
$ws.Range("B242").Value = "APPLE"
$ws.Range("C242").Value = "IPHONE XS"
$ws.Range("D242").Value = "64GB"

$ws.Range("A243").Value = "`niPhone XS Max 64 GB `n"
$ws.Range("B243").Value = "APPLE"
$ws.Range("C243").Value = "IPHONE XS MAX"
$ws.Range("D243").Value = "64GB"

$ws.Range("A244").Value = "`nCAT S52 Dual SIM`n"
$ws.Range("B244").Value = "CAT"
$ws.Range("C244").Value = "S52"
$ws.Range("D244").Value = "INFONA"

$ws.Range("A245").Value = "`nCAT S61 Dual SIM`n"
$ws.Range("B245").Value = "CAT"
$ws.Range("C245").Value = "S61"
$ws.Range("D245").Value = "INFONA"

$ws.Range("A246").Value = "`nDoro 632 Single SIM`n"
$ws.Range("B246").Value = "DORO"
$ws.Range("C246").Value = "'632"
$ws.Range("D246").Value = "INFONA"

$ws.Range("A247").Value = "`nHuawei P Smart Pro Dual SIM`n"
$ws.Range("B247").Value = "HUAWEI"
$ws.Range("C247").Value = "P SMART PRO"
$ws.Range("D247").Value = "INFONA"

$ws.Range("A248").Value = "`nHuawei  P Smart Z Dual SIM `n"
$ws.Range("B248").Value = "HUAWEI"
$ws.Range("C248").Value = "P SMART Z"
$ws.Range("D248").Value = "INFONA"

$ws.Range("A249").Value = "`nHuawei P30 Dual SIM`n"
$ws.Range("B249").Value = "HUAWEI"
$ws.Range("C249").Value = "P30"
$ws.Range("D249").Value = "INFONA"

$ws.Range("A250").Value = "`nHuawei P30 Lite Dual SIM `n"
$ws.Range("B250").Value = "HUAWEI"
$ws.Range("C250").Value = "P30 LITE"
$ws.Range("D250").Value = "INFONA"

$ws.Range("A251").Value = "`nHuawei P30 Pro Dual SIM`n"
$ws.Range("B251").Value = "HUAWEI"
$ws.Range("C251").Value = "P30 PRO"
$ws.Range("D251").Value = "INFONA"

$ws.Range("A252").Value = "`nHuawei Y6 2019 Dual SIM`n"
$ws.Range("B252").Value = "HUAWEI"
$ws.Range("C252").Value = "Y6"
$ws.Range("D252").Value = "INFONA"

$ws.Range("A253").Value = "`nHuawei Y7 2019 Dual SIM`n"
$ws.Range("B253").Value = "HUAWEI"
$ws.Range("C253").Value = "Y7"
$ws.Range("D253").Value = "INFONA"

$ws.Range("A254").Value = "`nSamsung Galaxy A10 Dual SIM`n"
$ws.Range("B254").Value = "SAMSUNG"
$ws.Range("C254").Value = "GALAXY A10"
$ws.Range("D254").Value = "INFONA"

$ws.Range("A255").Value = "`nSamsung Galaxy A30s Dual SIM`n"
$ws.Range("B255").Value = "SAMSUNG"
$ws.Range("C255").Value = "GALAXY A30S"
$ws.Range("D255").Value = "INFONA"

$ws.Range("A256").Value = "`nSamsung Galaxy A50 Dual SIM`n"
$ws.Range("B256").Value = "SAMSUNG"
$ws.Range("C256").Value = "GALAXY A50"
$ws.Range("D256").Value = "INFONA"

$ws.Range("A257").Value = "`nSamsung Galaxy A51 Dual SIM`n"
$ws.Range("B257").Value = "SAMSUNG"
$ws.Range("C257").Value = "GALAXY A51"
$ws.Range("D257").Value = "INFONA"

$ws.Range("A258").Value = "`nSamsung Galaxy A70 Dual SIM`n"
$ws.Range("B258").Value = "SAMSUNG"
$ws.Range("C258").Value = "GALAXY A70"
$ws.Range("D258").Value = "INFONA"

$ws.Range("A259").Value = "`nSamsung Galaxy A71`n"
$ws.Range("B259").Value = "SAMSUNG"
$ws.Range("C259").Value = "GALAXY A71"
$ws.Range("D259").Value = "INFONA"

$ws.Range("A260").Value = "`nSamsung Galaxy Fold`n"
$ws.Range("B260").Value = "SAMSUNG"
$ws.Range("C260").Value = "GALAXY FOLD"
$ws.Range("D260").Value = "INFONA"

$ws.Range("A261").Value = "`nSamsung Galaxy Note10 Dual SIM `n"
$ws.Range("B261").Value = "SAMSUNG"
$ws.Range("C261").Value = "GALAXY NOTE 10"
$ws.Range("D261").Value = "INFONA"

$ws.Range("A262").Value = "`nSamsung Galaxy Note 10 Lite Dual SIM`n"
$ws.Range("B262").Value = "SAMSUNG"
$ws.Range("C262").Value = "GALAXY NOTE 10 LITE"
$ws.Range("D262").Value = "INFONA"

$ws.Range("A263").Value = "`nSamsung Galaxy Note10+ Dual SIM`n"
$ws.Range("B263").Value = "SAMSUNG"
$ws.Range("C263").Value = "GALAXY NOTE 10+"
$ws.Range("D263").Value = "INFONA"

$ws.Range("A264").Value = "`nSamsung Galaxy S10 Dual SIM 128GB`n"
$ws.Range("B264").Value = "SAMSUNG"
$ws.Range("C264").Value = "GALAXY S10"
$ws.Range("D264").Value = "128GB"

$ws.Range("A265").Value = "`nSamsung Galaxy S10 Lite Dual SIM`n"
$ws.Range("B265").Value = "SAMSUNG"
$ws.Range("C265").Value = "GALAXY S10 LITE"
$ws.Range("D265").Value = "INFONA"

$ws.Range("A266").Value = "`nSamsung Galaxy S10+ Dual SIM 128 GB`n"
$ws.Range("B266").Value = "SAMSUNG"
$ws.Range("C266").Value = "GALAXY S10+"
$ws.Range("D266").Value = "128GB"

$ws.Range("A267").Value = "`nSamsung Galaxy S20`n"
$ws.Range("B267").Value = "SAMSUNG"
$ws.Range("C267").Value = "GALAXY S20"
$ws.Range("D267").Value = "INFONA"

$ws.Range("A268").Value = "`nSamsung Galaxy S20 Ultra 5G`n"
$ws.Range("B268").Value = "SAMSUNG"
$ws.Range("C268").Value = "GALAXY S20 ULTRA"
$ws.Range("D268").Value = "INFONA"

$ws.Range("A269").Value = "`nSamsung Galaxy S20+`n"
$ws.Range("B269").Value = "SAMSUNG"
$ws.Range("C269").Value = "GALAXY S20+"
$ws.Range("D269").Value = "INFONA"

$ws.Range("A270").Value = "`nSamsung Galaxy Z Flip`n"
$ws.Range("B270").Value = "SAMSUNG"
$ws.Range("C270").Value = "GALAXY Z FLIP"
$ws.Range("D270").Value = "INFONA"

$ws.Range("A271").Value = "`nXiaomi Mi Note 10 Dual SIM `n"
$ws.Range("B271").Value = "XIAOMI"
$ws.Range("C271").Value = "MI NOTE 10"
$ws.Range("D271").Value = "INFONA"

$ws.Range("A272").Value = "`nXiaomi Redmi 7A Dual SIM`n"
$ws.Range("B272").Value = "XIAOMI"
$ws.Range("C272").Value = "REDMI 7A"
$ws.Range("D272").Value = "INFONA"

$ws.Range("A273").Value = "`nXiaomi Redmi Note 8 Pro Dual SIM`n"
$ws.Range("B273").Value = "XIAOMI"
$ws.Range("C273").Value = "REDMI NOTE 8 PRO"
$ws.Range("D273").Value = "INFONA"

$ws.Range("A274").Value = "`nXiaomi Redmi Note 8T Dual SIM`n"
$ws.Range("B274").Value = "XIAOMI"
$ws.Range("C274").Value = "REDMI NOTE 8T"
$ws.Range("D274").Value = "INFONA"

$ws.Range("A275").Value = "`nMobitel Apple iPhone 11 128GB White "
$ws.Range("B275").Value = "APPLE"
$ws.Range("C275").Value = "IPHONE 11"
$ws.Range("D275").Value = "128GB"

$ws.Range("A276").Value = "`nMobitel Apple iPhone 11 64GB White "
$ws.Range("B276").Value = "APPLE"
$ws.Range("C276").Value = "IPHONE 11"
$ws.Range("D276").Value = "64GB"

$ws.Range("A277").Value = "`nMobitel Apple iPhone 8 Plus 64GB Space Grey "
$ws.Range("B277").Value = "APPLE"
$ws.Range("C277").Value = "IPHONE 8 PLUS"
$ws.Range("D277").Value = "64GB"

$ws.Range("A278").Value = "`nMobitel Samsung Galaxy A70 narančasti 128GB dual SIM SM-A705F "
$ws.Range("B278").Value = "SAMSUNG"
$ws.Range("C278").Value = "GALAXY A70"
$ws.Range("D278").Value = "128GB"

$ws.Range("A279").Value = "`n Mobitel Samsung Galaxy S10 kraljevsko crveni 128GB dual SIM SM-G973F "
$ws.Range("B279").Value = "SAMSUNG"
$ws.Range("C279").Value = "GALAXY S10"
$ws.Range("D279").Value = "128GB"

$ws.Range("A280").Value = "`nSamsung`n—Galaxy S20 DS sivi"
$ws.Range("B280").Value = "SAMSUNG"
$ws.Range("C280").Value = "GALAXY S20"
$ws.Range("D280").Value = "INFONA"

$ws.Range("A281").Value = "`nSamsung`n—Galaxy S20 Ultra DS sivi"
$ws.Range("B281").Value = "SAMSUNG"
$ws.Range("C281").Value = "GALAXY S20 ULTRA"
$ws.Range("D281").Value = "INFONA"

$ws.Range("A282").Value = "`nSamsung`n—Galaxy S20+ DS crni"
$ws.Range("B282").Value = "SAMSUNG"
$ws.Range("C282").Value = "GALAXY S20+"
$ws.Range("D282").Value = "INFONA"

$ws.Range("A283").Value = "`nMobitel XIAOMI MI 9T 6/64GB: CRNI "
$ws.Range("B283").Value = "XIAOMI"
$ws.Range("C283").Value = "MI 9T"
$ws.Range("D283").Value = "64GB"

$ws.Range("A284").Value = "`nMobitel Xiaomi Mi Note 10 6GB/128GB Zelena  "
$ws.Range("B284").Value = "XIAOMI"
$ws.Range("C284").Value = "MI NOTE 10"
$ws.Range("D284").Value = "128GB"
